$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "42.790.34"
$ws.Cells.Item(2,5).Value = "  -0.03%  "

$ws.Cells.Item(3,4).Value = "2.279.33"
$ws.Cells.Item(3,5).Value = "  +0.94%  "

$ws.Cells.Item(4,5).Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "251.28"
$ws.Cells.Item(5,5).Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "0.636"
$ws.Cells.Item(6,5).Value = "  +1.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "75.65"
$ws.Cells.Item(7,5).Value = "  +7.43%  "

$ws.Cells.Item(8,5).Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.647"
$ws.Cells.Item(9,5).Value = "  -1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "39.72"
$ws.Cells.Item(10,5).Value = "  +2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.0976"
$ws.Cells.Item(11,5).Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "7.39"
$ws.Cells.Item(12,5).Value = "  -1.01%  "

$ws.Cells.Item(13,5).Value = "  +1.27%  "

$ws.Cells.Item(14,4).Value = "2.623.60"
$ws.Cells.Item(14,5).Value = "  +1.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "15.10"
$ws.Cells.Item(15,5).Value = "  +1.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.870"
$ws.Cells.Item(16,5).Value = "  -1.01%  "

$ws.Cells.Item(17,4).Value = "2.273.31"
$ws.Cells.Item(17,5).Value = "  +0.68%  "

$ws.Cells.Item(18,4).Value = "42.707.36"
$ws.Cells.Item(18,5).Value = "  -0.07%  "

$ws.Cells.Item(19,4).Value = "0.0₃0997"
$ws.Cells.Item(19,5).Value = "  +0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "6.21"
$ws.Cells.Item(20,5).Value = "  -1.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "72.38"
$ws.Cells.Item(21,5).Value = "  -0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "236.48"
$ws.Cells.Item(22,5).Value = "  +0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "2.16"
$ws.Cells.Item(23,5).Value = "  +4.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "3.85"
$ws.Cells.Item(24,5).Value = "  -2.08%  "

$ws.Cells.Item(25,5).Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "11.28"
$ws.Cells.Item(26,5).Value = "  -1.41%  "

$ws.Cells.Item(27,5).Value = "  -1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "2.20"
$ws.Cells.Item(28,5).Value = "  +4.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "167.54"
$ws.Cells.Item(29,5).Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "21.02"
$ws.Cells.Item(30,5).Value = "  +0.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "0.0867"
$ws.Cells.Item(31,5).Value = "  +9.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "6.43"
$ws.Cells.Item(32,5).Value = "  -2.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "0.126"
$ws.Cells.Item(33,5).Value = "  -0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "31.97"
$ws.Cells.Item(34,5).Value = "  +2.08%  "

$ws.Cells.Item(35,5).Value = "  +1.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "4.54"
$ws.Cells.Item(36,5).Value = "  +2.58%  "

$ws.Cells.Item(37,5).Value = "  +1.01%  "

$ws.Cells.Item(38,5).Value = "  -4.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "13.53"
$ws.Cells.Item(39,5).Value = "  +8.95%  "

$ws.Cells.Item(40,5).Value = "  -0.65%  "

$ws.Cells.Item(41,5).Value = "  +1.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.207"
$ws.Cells.Item(42,5).Value = "  +2.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "61.64"
$ws.Cells.Item(43,5).Value = "  -0.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "8.93"
$ws.Cells.Item(44,5).Value = "  -1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "106.79"
$ws.Cells.Item(45,5).Value = "  +12.74%  "

$ws.Cells.Item(46,2).Value = "Cronos"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.101"
$ws.Cells.Item(46,5).Value = "  -1.34%  "

$ws.Cells.Item(47,2).Value = "FTXToken"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "4.64"
$ws.Cells.Item(47,5).Value = "  -4.13%  "

$ws.Cells.Item(48,5).Value = "  -0.39%  "

$ws.Cells.Item(49,5).Value = "  -0.58%  "

$ws.Cells.Item(50,5).Value = "  -1.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "4.21"
$ws.Cells.Item(51,5).Value = "  -1.99%  "
